# "Add files via upload" — the anagram quiz content was translated from
# Japanese to English (scrambled-word / answer pairs), and the active
# selection moved from B3 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "anagrams"
$ws.Range("B1").Value = "answer"

# Row 2: "plpae" is an anagram of "apple"
$ws.Range("A2").Value = "plpae"
$ws.Range("B2").Value = "apple"

# Row 3: "elnom" is an anagram of "melon"
$ws.Range("A3").Value = "elnom"
$ws.Range("B3").Value = "melon"

# Active cell moved to B2
$ws.Range("B2").Select()
